$d = $word.ActiveDocument

$d.Content.Find.Execute("Vector Machine", $true, $false, $false, $false, $false, $true, 1, $false, "Voting Classifier", 2)
